$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text cells (rich-text strings collapsed to plain text; content is what matters) ---
# A8 holds "Volume " & "31" & "   Number  " & "37" -> bump report number 37 -> 38
$ws.Range("A8").Value = "Volume 31   Number  38"
# C9 holds "Report Covering the Week  " & "9/9/2024" & "  Through  " & "9/15/2024" -> shift week by 7 days
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# --- Row 15 (Rape) ---
$ws.Range("M15").Value = -16.666666666666
$ws.Range("N15").Value = -68.085106382978

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 500
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 5.882352941176
$ws.Range("I16").Value = 169
$ws.Range("J16").Value = 139
$ws.Range("K16").Value = 21.582733812949
$ws.Range("L16").Value = 13.422818791946
$ws.Range("M16").Value = -0.588235294117
$ws.Range("N16").Value = -74.198473282442

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 18.75
$ws.Range("I17").Value = 336
$ws.Range("J17").Value = 293
$ws.Range("K17").Value = 14.675767918088
$ws.Range("L17").Value = 20.430107526881
$ws.Range("M17").Value = 111.320754716981
$ws.Range("N17").Value = -24.832214765100

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -27.777777777777
$ws.Range("I18").Value = 115
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 25
$ws.Range("L18").Value = -17.266187050359
$ws.Range("M18").Value = 61.971830985915
$ws.Range("N18").Value = -65.256797583081

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 18.181818181818
$ws.Range("F19").Value = 32
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 283
$ws.Range("J19").Value = 313
$ws.Range("K19").Value = -9.584664536741
$ws.Range("L19").Value = -19.373219373219
$ws.Range("M19").Value = 46.632124352331
$ws.Range("N19").Value = -33.096926713948

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 42.857142857142
$ws.Range("I20").Value = 62
$ws.Range("J20").Value = 61
$ws.Range("K20").Value = 1.639344262295
$ws.Range("L20").Value = -6.060606060606
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -79.054054054054

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 62.5
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = 4.716981132075
$ws.Range("I21").Value = 985
$ws.Range("J21").Value = 914
$ws.Range("K21").Value = 7.768052516411
$ws.Range("L21").Value = -1.696606786427
$ws.Range("M21").Value = 52.713178294573
$ws.Range("N21").Value = -55.670567056705

# --- Row 22 (Transit) --- C22 switches from the text placeholder "0" to a real numeric 1
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 20

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 30
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = -9.090909090909
$ws.Range("I23").Value = 288
$ws.Range("J23").Value = 286
$ws.Range("K23").Value = 0.699300699300
$ws.Range("L23").Value = -7.692307692307
$ws.Range("M23").Value = 48.453608247422

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = -36.231884057971
$ws.Range("I24").Value = 534
$ws.Range("J24").Value = 666
$ws.Range("K24").Value = -19.819819819819
$ws.Range("L24").Value = -16.431924882629
$ws.Range("M24").Value = 12.896405919661

# --- Row 25 (Retail Theft) ---
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -81.818181818181
$ws.Range("I25").Value = 101
$ws.Range("J25").Value = 202
$ws.Range("K25").Value = -50
$ws.Range("L25").Value = -35.256410256410

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 66
$ws.Range("H26").Value = -12.121212121212
$ws.Range("I26").Value = 567
$ws.Range("J26").Value = 456
$ws.Range("K26").Value = 24.342105263157
$ws.Range("L26").Value = 43.544303797468
$ws.Range("M26").Value = 11.394891944990

# --- Row 28 (Other Sex Crimes) --- D28/E28 switch to the text placeholders "0" / "***.*"
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 48
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = 2.127659574468

# --- Row 29 (Shooting Vic.) ---
$ws.Range("M29").Value = -62.068965517241

# --- Row 30 (Shooting Inc.) ---
$ws.Range("M30").Value = -69.230769230769
